# Applies the "per output weighting" / new experiment rows update to
# the "nn results 2020" workbook (sheet Foaie1 / sheet1.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foaie1")

# --- Header row tweaks -------------------------------------------------
# D1: "Weights (negative, positive)" -> "Weights"
$ws.Range("D1").Value = "Weights"

# M1: new "Notes" column header, styled like the rest of row 1 (center/center)
$ws.Range("A1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M1").Value = "Notes"

# --- Row 9: finish the previously-empty shared-formula row -------------
$ws.Range("A9").Value = 384
$ws.Range("B9").Value = 128
$ws.Range("C9").Value = "Conv, (128, 128), (64, 128), (1,128)"
$ws.Range("D9").Value = 4.62
$ws.Range("E9").Value = 0.25
$ws.Range("F9").Value = 0.886
$ws.Range("G9").Value = 0.7
$ws.Range("H9").Value = 0.58
$ws.Range("I9").Formula = "=2*(G9*H9)/(G9+H9)"
$ws.Range("J9").Value = "41,3k"
$ws.Range("L9").Value = 20

# --- Row 10 --------------------------------------------------------------
$ws.Range("A10").Value = 384
$ws.Range("B10").Value = 128
$ws.Range("C10").Value = "Conv, (128, 128), (64, 128), (1,128)"
$ws.Range("D10").Value = 4.62
$ws.Range("E10").Value = "?"
$ws.Range("F10").Value = 0.78
$ws.Range("G10").Value = 0.48
$ws.Range("H10").Value = 0.88
$ws.Range("I10").Formula = "=2*(G10*H10)/(G10+H10)"
$ws.Range("J10").Value = "41,3k"
$ws.Range("L10").Value = 20
$ws.Range("M10").Value = 'sample_weight_mode="temporal"'

# --- Row 11 ----------------------------------------------------------------
$ws.Range("A11").Value = 384
$ws.Range("B11").Value = 128
$ws.Range("C11").Value = "Conv, (128, 128), (64, 128), (1,128)"
$ws.Range("D11").Value = 8.28
$ws.Range("E11").Value = 0.2265
$ws.Range("F11").Value = 0.837
$ws.Range("G11").Value = 0.41
$ws.Range("H11").Value = 0.84
$ws.Range("I11").Formula = "=2*(G11*H11)/(G11+H11)"
$ws.Range("J11").Value = "41,3k"
$ws.Range("L11").Value = 20

# --- Row 12 ----------------------------------------------------------------
$ws.Range("A12").Value = 384
$ws.Range("B12").Value = 128
$ws.Range("C12").Value = "Conv, (128, 128), (64, 128), (1,128)"
$ws.Range("D12").Value = 4.86
$ws.Range("E12").Value = 0.3
$ws.Range("F12").Value = 0.88
$ws.Range("G12").Value = 0.66
$ws.Range("H12").Value = 0.59
$ws.Range("I12").Formula = "=2*(G12*H12)/(G12+H12)"
$ws.Range("J12").Value = "41,3k"
$ws.Range("L12").Value = 20
$ws.Range("M12").Value = "added random scaling of input"

# --- Row 13 ----------------------------------------------------------------
$ws.Range("A13").Value = 384
$ws.Range("B13").Value = 128
$ws.Range("C13").Value = "Conv, (128, 128), (64, 128), (32,128),(1,128), elu"
$ws.Range("D13").Value = 4.86
$ws.Range("E13").Value = 0.36
$ws.Range("F13").Value = 0.78
$ws.Range("G13").Value = 0.43
$ws.Range("H13").Value = 0.73
$ws.Range("I13").Formula = "=2*(G13*H13)/(G13+H13)"
$ws.Range("J13").Value = "43,1k"
$ws.Range("L13").Value = 20

# --- Row 14 ----------------------------------------------------------------
$ws.Range("A14").Value = 384
$ws.Range("B14").Value = 128
$ws.Range("C14").Value = "Conv, (96, 128), (32, 128), (1,128)"
$ws.Range("D14").Value = 3.38
$ws.Range("E14").Value = 0.35
$ws.Range("F14").Value = 0.79
$ws.Range("G14").Value = 0.55
$ws.Range("H14").Value = 0.79
$ws.Range("I14").Formula = "=2*(G14*H14)/(G14+H14)"
$ws.Range("J14").Value = "41,3k"
$ws.Range("L14").Value = 20
$ws.Range("M14").Value = "very reduced dataset"

# --- Cosmetics: column widths + selection ----------------------------------
# Column C grew (longer "elu" structure string), column D shrank back to the
# sheet default now that it only holds short "n, m" weight pairs.
$ws.Columns("C").ColumnWidth = 38
$ws.Columns("D").ColumnWidth = 8.5

$ws.Range("G23").Select()
